$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 259-260; this shifts all existing rows (259..276)
# down to (261..278) and extends the sheet dimension automatically.
$ws.Rows("259:260").Insert()

# Row 259: new weekly record (Modesto / Primera)
$ws.Cells.Item(259, 1).Value = 6
$ws.Cells.Item(259, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(259, 3).Value = "Metropolitana"
$ws.Cells.Item(259, 4).Value = 44946
$ws.Cells.Item(259, 5).Value = 13
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100103
$ws.Cells.Item(259, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(259, 9).Value = 100103003
$ws.Cells.Item(259, 10).Value = "Damasco"
$ws.Cells.Item(259, 11).Value = "Modesto"
$ws.Cells.Item(259, 12).Value = "Primera"
$ws.Cells.Item(259, 13).Value = 290
$ws.Cells.Item(259, 14).Value = 14000
$ws.Cells.Item(259, 15).Value = 15000
$ws.Cells.Item(259, 16).Value = 14500
$ws.Cells.Item(259, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(259, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(259, 19).Value = 906
$ws.Cells.Item(259, 20).Value = 16

# Row 260: new weekly record (Modesto / Segunda)
$ws.Cells.Item(260, 1).Value = 6
$ws.Cells.Item(260, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(260, 3).Value = "Metropolitana"
$ws.Cells.Item(260, 4).Value = 44946
$ws.Cells.Item(260, 5).Value = 13
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100103
$ws.Cells.Item(260, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(260, 9).Value = 100103003
$ws.Cells.Item(260, 10).Value = "Damasco"
$ws.Cells.Item(260, 11).Value = "Modesto"
$ws.Cells.Item(260, 12).Value = "Segunda"
$ws.Cells.Item(260, 13).Value = 275
$ws.Cells.Item(260, 14).Value = 11000
$ws.Cells.Item(260, 15).Value = 11000
$ws.Cells.Item(260, 16).Value = 11000
$ws.Cells.Item(260, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(260, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(260, 19).Value = 688
$ws.Cells.Item(260, 20).Value = 16
